# Jeannine's log - add TUESDAY (Sept 27, 2016 / serial 42640) entries
# mirrors the day-separator + data-row pattern already used throughout the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- 1. Build the new rows' formatting by cloning existing template rows ---
# Row 5 is a "day separator" row (empty cells + day name in col B)
$ws.Range("A5:F5").Copy()
$ws.Range("A259:F259").PasteSpecial(-4122)

# Rows 256-269 (ten data rows) get the same formatting as the last existing
# data row (255)
$ws.Range("A255:F255").Copy()
$ws.Range("A260:F269").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 2. Day separator row ---
$ws.Range("B259").Value = "TUESDAY"

# --- 3. Data rows (Task Type, Date, Time, Building, Room, Comments) ---
$ws.Range("A260").Value = "Demo"
$ws.Range("B260").Value = 42640
$ws.Range("C260").Value = "1630"
$ws.Range("D260").Value = "CC"
$ws.Range("E260").Value = "109"
$ws.Range("F260").Value = "Please demo equipment to client and make sure client is happy."

$ws.Range("A261").Value = "AV Shutdown"
$ws.Range("B261").Value = 42640
$ws.Range("C261").Value = "1630"
$ws.Range("D261").Value = "LSB"
$ws.Range("E261").Value = "107"
$ws.Range("F261").Value = "Make sure neck mic goes back to drawer and log off touchscreen."

$ws.Range("A262").Value = "AV Shutdown"
$ws.Range("B262").Value = 42640
$ws.Range("C262").Value = "1730"
$ws.Range("D262").Value = "LSB"
$ws.Range("E262").Value = "106"
$ws.Range("F262").Value = "Make sure neck mic goes back to drawer and log off touchscreen."

$ws.Range("A263").Value = "AV Shutdown"
$ws.Range("B263").Value = 42640
$ws.Range("C263").Value = "1900"
$ws.Range("D263").Value = "LSB"
$ws.Range("E263").Value = "103"
$ws.Range("F263").Value = "Make sure neck mic goes back to drawer and log off touchscreen."

$ws.Range("A264").Value = "Demo"
$ws.Range("B264").Value = 42640
$ws.Range("C264").Value = "1800"
$ws.Range("D264").Value = "CLH"
$ws.Range("E264").Value = "M"
$ws.Range("F264").Value = "NO CEILING PROJECTOR - Use roll in PC and Projector that is in room. Make sure client is okay."

$ws.Range("A265").Value = "Demo"
$ws.Range("B265").Value = 42640
$ws.Range("C265").Value = "1900"
$ws.Range("D265").Value = "CLH"
$ws.Range("E265").Value = "J"
$ws.Range("F265").Value = "NO CEILING PROJECTOR - Use roll in PC and Projector that is in room. Make sure client is okay."

$ws.Range("A266").Value = "Lockup"
$ws.Range("B266").Value = 42640
$ws.Range("C266").Value = "2000"
$ws.Range("D266").Value = "CLH"
$ws.Range("E266").Value = "K"
$ws.Range("F266").Value = "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS."

$ws.Range("A267").Value = "Lockup"
$ws.Range("B267").Value = 42640
$ws.Range("C267").Value = "2130"
$ws.Range("D267").Value = "CLH"
$ws.Range("E267").Value = "M"
$ws.Range("F267").Value = "LEAVE ROLL IN PC AND PROJECTOR IN ROOM - JUST TURN OFF. PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lasonde 1011 office. PLEASE LOCK ALL 4 DOORS."

$ws.Range("A268").Value = "Lockup"
$ws.Range("B268").Value = 42640
$ws.Range("C268").Value = "2200"
$ws.Range("D268").Value = "CLH"
$ws.Range("E268").Value = "J"
$ws.Range("F268").Value = "LEAVE ROLL IN PC AND PROJECTOR IN ROOM - JUST TURN OFF. PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lasonde 1011 office. PLEASE LOCK ALL 4 DOORS."

$ws.Range("A269").Value = "Lockup"
$ws.Range("B269").Value = 42640
$ws.Range("C269").Value = "2200"
$ws.Range("D269").Value = "CLH"
$ws.Range("E269").Value = "H"
$ws.Range("F269").Value = "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS."

# --- 4. Wrapped comments push a couple of rows taller, same as the rest of
#        the sheet (rows auto-grow with wrapped text; set explicitly here to
#        mirror the authored heights) ---
$ws.Range("A264:F264").EntireRow.RowHeight = 30
$ws.Range("A265:F265").EntireRow.RowHeight = 30
$ws.Range("A266:F266").EntireRow.RowHeight = 30
$ws.Range("A267:F267").EntireRow.RowHeight = 45
$ws.Range("A268:F268").EntireRow.RowHeight = 45
$ws.Range("A269:F269").EntireRow.RowHeight = 30

# --- 5. Update the view: scroll down and move the active selection, as in
#        the authored workbook (topLeftCell A254, active cell F272) ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 254
$win.ScrollColumn = 1
$ws.Range("F272").Select()
